# Rewrite the sales report rows 2-10 with the new transaction data,
# keeping the new order-date / amount / discount columns stored as TEXT
# (matching the original workbook's convention of storing these as
# shared-string text rather than numeric cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Force "text" number format on the cells whose values look like
#        numbers or dates (Order Date / Amount / Discount columns for the
#        data rows), so Excel doesn't silently coerce them into numeric or
#        date serial values when we assign the string below.
#        NOTE: multi-area ranges ("D2:D10,G2:G10,H2:H10") only apply bulk
#        property writes to the first area here, so each column is handled
#        as its own contiguous range.
$colD = $ws.Range("D2:D10")
$colG = $ws.Range("G2:G10")
$colH = $ws.Range("H2:H10")
$colD.NumberFormat = "@"
$colG.NumberFormat = "@"
$colH.NumberFormat = "@"

# ---- 2. New transaction rows (SL No, Order ID, Customer, Order Date,
#         Product, Qty, Amount, Discount, Payment Method, Status)
$data = @(
    @(1, "#FBZ-115EG1Y", "Joyal K",        "2025-04-06", "Kettle Bell 6 Kgs for home gym", 1, "1,008", "150", "netbanking", "Pending"),
    @(2, "#FBZ-117E3WP", "Joyal K",        "2025-04-05", "Sports Premium Blitz Blue Kettlebell 8kg Cast Iron Vinyl Coated Solid Kettlebell", 1, "1,635", "150", "netbanking", "Delivered"),
    @(3, "#FBZ-PBMIZYX", "Abhinav K",      "2025-04-04", "Cast Iron Vinyl Coated Dumbbells for gym Workout", 1, "599", "200", "wallet", "Pending"),
    @(4, "#FBZ-TM3CGR9", "Cezanne P",      "2025-04-03", "MEDIX Soft Medicine Ball (2), Rubber for Adults", 1, "629", "200", "netbanking", "Pending"),
    @(5, "#FBZ-VDGHF7G", "Achyuth J",      "2025-04-02", "Rubber Medicine Ball Weights for men & women", 1, "2,052", "125", "cod", "Pending"),
    @(6, "#FBZ-XIPA8OQ", "Achu K",         "2025-04-01", "Pair of two PVC Dumbbells Set Hex for all", 1, "1,065", "0", "netbanking", "Pending"),
    @(7, "#FBZ-9VZ5PY4", "Basim M",        "2025-04-07", "Holistic Fitness 6kg Rubber Slam Ball", 2, "8,698", "500", "cod", "Delivered"),
    @(8, "#FBZ-JLTUPXB", "Joyal Kuriakose","2025-04-07", "Adjustable Hand Grip Strengthener, Hand Gripper With Counter for Men", 3, "607", "200", "netbanking", "Cancelled"),
    @(9, "#FBZ-L0SXRRL", "Joyal Kuriakose","2025-04-07", "Adjustable Hand Grip Strengthener, Hand Gripper With Counter for Men", 1, "295", "0", "netbanking", "Delivered")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]   # A - SL No
    $ws.Cells.Item($row, 2).Value = $rec[1]   # B - Order ID
    $ws.Cells.Item($row, 3).Value = $rec[2]   # C - Customer
    $ws.Cells.Item($row, 4).Value = $rec[3]   # D - Order Date
    $ws.Cells.Item($row, 5).Value = $rec[4]   # E - Product
    $ws.Cells.Item($row, 6).Value = $rec[5]   # F - Qty
    $ws.Cells.Item($row, 7).Value = $rec[6]   # G - Amount (Rs)
    $ws.Cells.Item($row, 8).Value = $rec[7]   # H - Discount (Rs)
    $ws.Cells.Item($row, 9).Value = $rec[8]   # I - Payment Method
    $ws.Cells.Item($row, 10).Value = $rec[9]  # J - Status
    $row = $row + 1
}

# ---- 3. Clear the explicit "text" style we applied above so the cells
#         fall back to the default (unstyled) cell format, same as the
#         rest of the sheet - only the underlying stored type stays text.
$colD.Style = "Normal"
$colG.Style = "Normal"
$colH.Style = "Normal"

# ---- 4. Updated summary block (rows 12-15): total orders stays 9, but the
#         amount/discount/net totals reflect the new data set.
$ws.Range("G13").Value = "₹16,588"
$ws.Range("G14").Value = "₹1,525"
$ws.Range("G15").Value = "₹15,063"
